$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found for replace:" $oldText
    }
}

function Insert-Paragraph-After($anchorText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: anchor text not found for insert:" $anchorText
        return
    }
    $para = $rng.Paragraphs(1)
    $para.Range.InsertParagraphAfter()
    $nextPara = $para.Next()
    $nextPara.Range.Text = $newText
}

function Delete-Paragraph($text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: text not found for delete:" $text
        return
    }
    $para = $rng.Paragraphs(1)
    $para.Range.Delete()
}

# ---------------------------------------------------------------------------
# Staff section (numId 1004)
# ---------------------------------------------------------------------------
Replace-Text "Meag Doherty, Chief of Staff" "Aaron Ogle, Director of Product"
Insert-Paragraph-After "Aaron Ogle, Director of Product" "Patrick Bateman, Deputy Executive Director & Chief of Staff"
Replace-Text "Mary Kate Mezzetti, Intern" "Kunal Kothari, Fellow"

# ---------------------------------------------------------------------------
# Board of Directors section (numId 1005)
# ---------------------------------------------------------------------------
Replace-Text "Darrell Issa, Chairman" "Phaedra Chrousos, Chairman"
Replace-Text "Seamus Kraft, President and Vice-Chairman" "Laurent Crenshaw, Vice Chairman"
Insert-Paragraph-After "Laurent Crenshaw, Vice Chairman" "Seamus Kraft, Secretary"
Replace-Text "James Lacy, Counsel" "Jo-Marie St. Martin, Member"
Delete-Paragraph "Tom Davis, Member"
Delete-Paragraph "Abhi Nemani, Member"
Delete-Paragraph "Larry Brady, Member"

# ---------------------------------------------------------------------------
# Board of Advisers section (numId 1006)
# ---------------------------------------------------------------------------
# Delete/remove entries first while text is still unambiguous, to avoid
# collisions with newly-introduced duplicate text from later replacements.
Delete-Paragraph "Aaron Bartnick"
Delete-Paragraph "Lanham Napier, BuildGroup"
Delete-Paragraph "Dr. Anne Washington, George Mason University"

Replace-Text "Dr. Anne Washington, Legal Data and Informatics" "Karien Bezuidenhout, Internal Growth"
Replace-Text "Joe Trippi, Political Campaigns" "Dr. Anne Washington, Open Legal Data"
Replace-Text "Karien Bezuidenhout, Shuttleworth Foundation" "Aaron Bartnick"
Replace-Text "Brandon Andrews, Values Partnership" "Brandon Andrews"
